# Updates the "cryptos" sheet with refreshed price/volume figures
# (and a couple of row re-orderings) per the Tue Oct 31 12:57:06 UTC 2023
# GitHub Actions data refresh.
#
# Note: column D ("Price") values are stored as text (e.g. "1.00", "246.00",
# "34.463.01") even though many look numeric. Plain `.Value = "..."` on a
# numeric-looking string causes Excel to silently coerce it to a real number
# (dropping trailing zeros / mis-parsing thousands-dot grouping), so for
# those cells we force the Text number format first, assign the string,
# then restore the cell's normal (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.463.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.806.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.601"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.44%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.298"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0698"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0966"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.068.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.824.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.649"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.468.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0790"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0525"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.388.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.662"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0188"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "82.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.954"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0499"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.968.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0127"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.38%  "
